$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split " ... pesquisar o trending top de palavras mais faladas ..."
#    so the English loanword "trending" sits in its own run (Word marks it
#    as a foreign/unknown word during spellcheck -> separate run boundary).
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("trending")
$r.Bold = 1
$r.Bold = 0

# ---------------------------------------------------------------------------
# 2) Split "Logar e cadastrar" into "Logar" + " e cadastrar" runs.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Logar")
$r2.Bold = 1
$r2.Bold = 0

# ---------------------------------------------------------------------------
# 3) Resize the traceability-matrix table (3rd table in the document):
#    table width/indent + last grid column + every cell width in every row.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(3)

$t.PreferredWidthType = 3
$t.PreferredWidth = 418.15
$t.Rows.LeftIndent = 13.95

$widths = @(58.25, 57.5, 57.5, 57.5, 57.55, 57.55, 72.3)
for ($i = 1; $i -le 7; $i++) {
  $col = $t.Columns.Item($i)
  $col.PreferredWidthType = 3
  $col.PreferredWidth = $widths[$i - 1]
  $col.Cells.PreferredWidthType = 3
  $col.Cells.PreferredWidth = $widths[$i - 1]
}
